$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "Dahi Bada"
$ws.Range("A9").Value = "Masala Dosa"
$ws.Range("A12").Value = "Pulav"

$ws.Activate()
$ws.Range("A12").Select()
